$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "Förändrad" (C) column for all existing data rows (2-244)
#    from 45190 to 45192.
for ($r = 2; $r -le 244; $r++) {
    $ws.Cells.Item($r, 3).Value = 45192
}

# 2) Row 244 gains an explicit row height (ht="15" customHeight="1").
$ws.Rows.Item(244).RowHeight = 15

# 3) Append three new notification rows (245-247).

# Row 245
$ws.Rows.Item(245).RowHeight = 15
$ws.Range("A245").Value = "A 45065-2023"
$ws.Range("B245").Value = 45191
$ws.Range("B245").NumberFormat = "YYYY-MM-DD"
$ws.Range("C245").Value = 45192
$ws.Range("C245").NumberFormat = "YYYY-MM-DD"
$ws.Range("D245").Value = "VÄSTMANLANDS LÄN"
$ws.Range("E245").Value = "NORBERG"
$ws.Range("F245").Value = "Bergvik skog väst AB"
$ws.Range("G245").Value = 4.5
$ws.Range("H245").Value = 0
$ws.Range("I245").Value = 0
$ws.Range("J245").Value = 0
$ws.Range("K245").Value = 0
$ws.Range("L245").Value = 0
$ws.Range("M245").Value = 0
$ws.Range("N245").Value = 0
$ws.Range("O245").Value = 0
$ws.Range("P245").Value = 0
$ws.Range("Q245").Value = 0
$ws.Range("R245").WrapText = $true

# Row 246
$ws.Rows.Item(246).RowHeight = 15
$ws.Range("A246").Value = "A 45051-2023"
$ws.Range("B246").Value = 45191
$ws.Range("B246").NumberFormat = "YYYY-MM-DD"
$ws.Range("C246").Value = 45192
$ws.Range("C246").NumberFormat = "YYYY-MM-DD"
$ws.Range("D246").Value = "VÄSTMANLANDS LÄN"
$ws.Range("E246").Value = "NORBERG"
$ws.Range("F246").Value = "Bergvik skog väst AB"
$ws.Range("G246").Value = 8.9
$ws.Range("H246").Value = 0
$ws.Range("I246").Value = 0
$ws.Range("J246").Value = 0
$ws.Range("K246").Value = 0
$ws.Range("L246").Value = 0
$ws.Range("M246").Value = 0
$ws.Range("N246").Value = 0
$ws.Range("O246").Value = 0
$ws.Range("P246").Value = 0
$ws.Range("Q246").Value = 0
$ws.Range("R246").WrapText = $true

# Row 247 (last row - no explicit custom row height, matching source pattern
# where only the final row of the sheet lacks ht="15" customHeight="1")
$ws.Range("A247").Value = "A 45066-2023"
$ws.Range("B247").Value = 45191
$ws.Range("B247").NumberFormat = "YYYY-MM-DD"
$ws.Range("C247").Value = 45192
$ws.Range("C247").NumberFormat = "YYYY-MM-DD"
$ws.Range("D247").Value = "VÄSTMANLANDS LÄN"
$ws.Range("E247").Value = "NORBERG"
$ws.Range("F247").Value = "Bergvik skog väst AB"
$ws.Range("G247").Value = 3.7
$ws.Range("H247").Value = 0
$ws.Range("I247").Value = 0
$ws.Range("J247").Value = 0
$ws.Range("K247").Value = 0
$ws.Range("L247").Value = 0
$ws.Range("M247").Value = 0
$ws.Range("N247").Value = 0
$ws.Range("O247").Value = 0
$ws.Range("P247").Value = 0
$ws.Range("Q247").Value = 0
$ws.Range("R247").WrapText = $true
